$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows appended to the bottom of the table (date 2025-10-22 / serial 45952)
$row104 = 683.76599999999996, 1435.3520000000003, 421.40200000000004, 262.94200000000001, 338.40600000000001, 618.05100000000004, 651.89700000000005, 221.23400000000001, 192.49299999999999, 297.64599999999996, 211.80999999999997, 188.45999999999998, 795.20499999999993, 1107.029, 628.84599999999989, 255.65899999999999, 179.59799999999998, 189.03199999999995, 69.89, 63.78, 29.27, 74.259999999999991, 59.76, 1.48
$row105 = 311.86599999999999, 259.83999999999997, 113.44499999999999, 158.54599999999999, 80.239999999999995, 232.22499999999999, 317.279, 193.02, 191.702, 204.23100000000002, 156.04400000000001, 372.18600000000004, 417.92300000000006, 832.17100000000005, 290.19600000000003, 295.27899999999994, 313.53699999999992, 182.071, 16.414000000000001, 65.305000000000007, 72.326999999999998, 26.658999999999999, 36.389000000000003, 36.136000000000003

$dateValue = Get-Date -Year 2025 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0

# Row 104: 四方坪站充电量(kw)
$ws.Cells.Item(104, 1).Value = $dateValue
$ws.Cells.Item(104, 1).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(104, 2).Value = "四方坪站充电量(kw)"
for ($i = 0; $i -lt $row104.Length; $i++) {
    $ws.Cells.Item(104, 3 + $i).Value = $row104[$i]
}

# Row 105: 高岭站充电量(kw)
$ws.Cells.Item(105, 1).Value = $dateValue
$ws.Cells.Item(105, 1).NumberFormat = "yyyy\-mm\-dd"
$ws.Cells.Item(105, 2).Value = "高岭站充电量(kw)"
for ($i = 0; $i -lt $row105.Length; $i++) {
    $ws.Cells.Item(105, 3 + $i).Value = $row105[$i]
}

# Update the view so scrolled position / selection mirror the authored workbook
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 71
$ws.Range("G109").Select()
